$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.751.34"
$ws.Range("E2").Value = "  +0.94%  "
$ws.Range("D3").Value = "1.702.03"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3933"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4053"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.38%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.519"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.005"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.70"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.67%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08863"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.62%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.473"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.70"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.101"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001324"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.86%  "
$ws.Range("D17").Value = "1.698.99"
$ws.Range("E17").Value = "  +0.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "99.49"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07066"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.81"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.083"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.003"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.77"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.62%  "
$ws.Range("D24").Value = "24.750.41"
$ws.Range("E24").Value = "  +0.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.165"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.62%  "
$ws.Range("E26").Value = "  +1.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.91%  "
$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.262"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +24.16%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "164.58"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "135.77"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.160"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.933"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09058"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.63%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.071"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02981"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2772"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.61%  "
$ws.Range("E37").Value = "  +0.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "11.06"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.47"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09253"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.466"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7749"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.24"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7212"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.600"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.214"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.46%  "
$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.001"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("B48").Value = "Flow"
$ws.Range("C48").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.350"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "140.22"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07989"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "89.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.34%  "
